$d = $word.ActiveDocument

# The original paragraph holds a single run whose <w:t> illegally nests
# <w:br/> elements (and literal marker text) inside the text node. The
# edit "consolidates non-text run tags": the leading plain-text span
# ("ellohay ") and the trailing plain-text span (" orldway.") are split
# out into their own plain runs, leaving the <w:br/>/marker-bearing
# middle section in its own (bold) run.
#
# Because that middle fragment is not representable through the
# Range.Text/Find object model (it contains non-text child nodes inside
# <w:t>), we rebuild the paragraph's raw WordprocessingML directly and
# push it back in with Range.InsertXML, which replaces the contents of
# the target range with exactly the XML supplied.

$para = $d.Paragraphs.Item(1)
$rng = $para.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00550C5A" w:rsidRDefault="0081227A"><w:r><w:t xml:space="preserve">ellohay </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t><w:br/>[MARKER_OPENING 0:&lt;w:r w:rsidRPr="0081227A">&lt;w:rPr>&lt;w:b/>&lt;w:bCs/>&lt;/w:rPr>&lt;w:t>]<w:br/>oldbay<w:br/>[MARKER_CLOSING 1:&lt;/w:t>&lt;/w:r>]<w:br/></w:t></w:r><w:r><w:t xml:space="preserve"> orldway.</w:t></w:r></w:p>'

[void]$rng.InsertXML($xml)
